$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 7370
$ws.Range("I32").Value = 7132.1665
$ws.Range("J32").Value = 7499.727
$ws.Range("K32").Value = 7132.1665
$ws.Range("L32").Value = 7499.727
$ws.Range("M32").Value = -6806.1665
$ws.Range("N32").Value = -8151.727

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2748.75
$ws.Range("I106").Value = 2572.3809
$ws.Range("J106").Value = 3983.3333
$ws.Range("K106").Value = 2572.3809
$ws.Range("L106").Value = 3983.3333
$ws.Range("M106").Value = -1941.3809
$ws.Range("N106").Value = -5245.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 9940.370000000001
$ws.Range("I116").Value = 8936.9375
$ws.Range("J116").Value = 11399.909
$ws.Range("K116").Value = 8936.9375
$ws.Range("L116").Value = 11399.909
$ws.Range("M116").Value = -5494.9375
$ws.Range("N116").Value = -18283.909

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 244385
$ws.Range("J134").Value = 244385
$ws.Range("L134").Value = 244385
$ws.Range("N134").Value = -254525

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2512.6316
$ws.Range("I137").Value = 1843.875
$ws.Range("K137").Value = 5531.625
$ws.Range("M137").Value = -2981.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 79195.47
$ws.Range("I32").Value = 84163.48
$ws.Range("K32").Value = 84163.48
$ws.Range("M32").Value = -83876.48

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1910.4445
$ws.Range("I110").Value = 1916
$ws.Range("K110").Value = 1916
$ws.Range("M110").Value = 129

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 234615.39
$ws.Range("J134").Value = 234615.39
$ws.Range("L134").Value = 234615.39
$ws.Range("N134").Value = -244755.39

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7275.4736
$ws.Range("I20").Value = 7784.8237
$ws.Range("J20").Value = 2946
$ws.Range("K20").Value = 7784.8237
$ws.Range("L20").Value = 2946
$ws.Range("M20").Value = -7537.8237
$ws.Range("N20").Value = -3440

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 27062.375
$ws.Range("J100").Value = 27062.375
$ws.Range("L100").Value = 27062.375
$ws.Range("N100").Value = -29226.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4618.1816
$ws.Range("I107").Value = 3718.5186
$ws.Range("J107").Value = 8666.666999999999
$ws.Range("K107").Value = 3718.5186
$ws.Range("L107").Value = 8666.666999999999
$ws.Range("M107").Value = -1798.5186
$ws.Range("N107").Value = -12506.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 27528.6
$ws.Range("J28").Value = 27528.6
$ws.Range("L28").Value = 27528.6
$ws.Range("N28").Value = -28018.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4222.4165
$ws.Range("I31").Value = 2303.1875
$ws.Range("K31").Value = 2303.1875
$ws.Range("M31").Value = -2008.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4222.4165
$ws.Range("I34").Value = 2303.1875
$ws.Range("K34").Value = 2303.1875
$ws.Range("M34").Value = -2101.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2674.4075
$ws.Range("I58").Value = 2635.1738
$ws.Range("J58").Value = 2900
$ws.Range("K58").Value = 2635.1738
$ws.Range("L58").Value = 2900
$ws.Range("M58").Value = -2432.1738
$ws.Range("N58").Value = -3306

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1251.8889
$ws.Range("I94").Value = 1070.3334
$ws.Range("J94").Value = 1342.6666
$ws.Range("K94").Value = 1070.3334
$ws.Range("L94").Value = 1342.6666
$ws.Range("M94").Value = -619.3334
$ws.Range("N94").Value = -2244.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2674.4075
$ws.Range("I136").Value = 2635.1738
$ws.Range("J136").Value = 2900
$ws.Range("K136").Value = 7905.5214
$ws.Range("L136").Value = 8700
$ws.Range("M136").Value = -5355.5214
$ws.Range("N136").Value = -13800

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1857.1428
$ws.Range("J80").Value = 1866.6666
$ws.Range("L80").Value = 5599.9998
$ws.Range("N80").Value = -7471.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 1857.1428
$ws.Range("J83").Value = 1866.6666
$ws.Range("L83").Value = 16799.9994
$ws.Range("N83").Value = -26159.9994

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 27756.25
$ws.Range("J95").Value = 29507.285
$ws.Range("L95").Value = 88521.855
$ws.Range("N95").Value = -92639.855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5652.769
$ws.Range("I131").Value = 1059.9
$ws.Range("J131").Value = 8523.3125
$ws.Range("K131").Value = 3179.7
$ws.Range("L131").Value = 25569.9375
$ws.Range("M131").Value = 1860.3
$ws.Range("N131").Value = -35649.9375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 14550
$ws.Range("J92").Value = 14550
$ws.Range("L92").Value = 14550
$ws.Range("N92").Value = -18294

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 65700
$ws.Range("J101").Value = 65700
$ws.Range("L101").Value = 65700
$ws.Range("N101").Value = -72190

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 2750.1428
$ws.Range("J13").Value = 6750
$ws.Range("L13").Value = 6750
$ws.Range("N13").Value = -7030

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4064.923
$ws.Range("I22").Value = 1587
$ws.Range("J22").Value = 5166.222
$ws.Range("K22").Value = 1587
$ws.Range("L22").Value = 5166.222
$ws.Range("M22").Value = -1292
$ws.Range("N22").Value = -5756.222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 4064.923
$ws.Range("I27").Value = 1587
$ws.Range("J27").Value = 5166.222
$ws.Range("K27").Value = 1587
$ws.Range("L27").Value = 5166.222
$ws.Range("M27").Value = -1480
$ws.Range("N27").Value = -5380.222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2205.7307
$ws.Range("I46").Value = 1883.3334
$ws.Range("J46").Value = 2482.0715
$ws.Range("K46").Value = 1883.3334
$ws.Range("L46").Value = 2482.0715
$ws.Range("M46").Value = -1695.3334
$ws.Range("N46").Value = -2858.0715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 17000
$ws.Range("J97").Value = 17000
$ws.Range("L97").Value = 17000
$ws.Range("N97").Value = -18982

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 22639.242
$ws.Range("J136").Value = 65520.5
$ws.Range("L136").Value = 196561.5
$ws.Range("N136").Value = -201661.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 29995
$ws.Range("J15").Value = 29995
$ws.Range("L15").Value = 29995
$ws.Range("N15").Value = -30571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 15785.429
$ws.Range("J54").Value = 26570.857
$ws.Range("L54").Value = 26570.857
$ws.Range("N54").Value = -27610.857

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 24000
$ws.Range("I61").Value = 24000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 24000
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -23708

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3135.5334
$ws.Range("I81").Value = 2616.6365
$ws.Range("J81").Value = 4562.5
$ws.Range("K81").Value = 5233.273
$ws.Range("L81").Value = 9125
$ws.Range("M81").Value = -4172.273
$ws.Range("N81").Value = -11247

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 125026936
$ws.Range("J82").Value = 125026936
$ws.Range("L82").Value = 125026936
$ws.Range("N82").Value = -125027702

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3135.5334
$ws.Range("I84").Value = 2616.6365
$ws.Range("J84").Value = 4562.5
$ws.Range("K84").Value = 26166.365
$ws.Range("L84").Value = 45625
$ws.Range("M84").Value = -20862.365
$ws.Range("N84").Value = -56233

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H85").Value = 125026936
$ws.Range("J85").Value = 125026936
$ws.Range("L85").Value = 125026936
$ws.Range("N85").Value = -125029588

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4142.263
$ws.Range("I96").Value = 2646.889
$ws.Range("K96").Value = 2646.889
$ws.Range("M96").Value = -1273.889

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 57399.8
$ws.Range("I125").Value = 38999.5
$ws.Range("J125").Value = 69666.664
$ws.Range("K125").Value = 38999.5
$ws.Range("L125").Value = 69666.664
$ws.Range("M125").Value = -34079.5
$ws.Range("N125").Value = -79506.664
